$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing label cell (A2) to the new label cells (A3:A5)
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update existing row 2 values
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 31

# Row 3 (new data row)
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 29

# Row 4 (new data row, previously row 3's value shifted down)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 24

# Row 5 (new data row)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 15
